$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("2€")

# Update "2€" column (G) flag values from 0 to 1 for specific rows
$ws.Range("G3").Value = 1
$ws.Range("G9").Value = 1
$ws.Range("G12").Value = 1
$ws.Range("G13").Value = 1
$ws.Range("G14").Value = 1
$ws.Range("G15").Value = 1

# Update the active selection on the sheet to H16 in the bottom-right frozen pane
$ws.Range("H16").Select()
